# Leave Card update - 12/22/2023 10:59 AM
# Shifts the recurring monthly PERIOD dates (col A, rows 94-107) forward to the
# end-of-month date, fills in the EARNED/Undertime figures for the newly
# completed Sep/Oct/Nov 2023 periods (rows 103-105), and records the Nov 2023
# "FL(1-0-0)" leave particular with its remarks date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Column A: PERIOD dates rolled forward one month (to month-end) ---
$ws.Range("A94").Value  = 44957
$ws.Range("A95").Value  = 44985
$ws.Range("A96").Value  = 45016
$ws.Range("A97").Value  = 45046
$ws.Range("A99").Value  = 45077
$ws.Range("A100").Value = 45107
$ws.Range("A101").Value = 45138
$ws.Range("A102").Value = 45169
$ws.Range("A103").Value = 45199
$ws.Range("A104").Value = 45230
$ws.Range("A105").Value = 45260
$ws.Range("A106").Value = 45291
$ws.Range("A107").Value = 45322
# (rows 106 / 107 dates set above)

# --- Row 103: EARNED posted for the period (EARNED 2 column recalculates) ---
$ws.Range("C103").Value = 1.25

# --- Row 104: EARNED posted for the period ---
$ws.Range("C104").Value = 1.25

# --- Row 105: new leave entry "FL(1-0-0)" with EARNED + undertime + remarks date ---
$ws.Range("B105").Value = "FL(1-0-0)"
$ws.Range("C105").Value = 1.25
$ws.Range("D105").Value = 1

# K105 picks up the date number format used by the other REMARKS-date cells
# in this table (e.g. K103) before being given its value.
$ws.Range("K103").Copy()
$ws.Range("K105").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("K105").Value = 45258

# Restore the on-screen selection to match where the edit left off.
$ws.Range("I9").Select()
$ws.Range("F113").Select()
